$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cell -> new literal text value, taken from the authoritative diff
# of the refreshed "cryptos" price/volume snapshot (coinranking.com scrape).
$updates = @(
    @('D2', '27.396.71'),
    @('E2', '  +1.98%  '),
    @('D3', '1.825.89'),
    @('E3', '  +1.04%  '),
    @('D4', '0.9999'),
    @('E4', '  -0.07%  '),
    @('D5', '312.79'),
    @('E5', '  +0.82%  '),
    @('D6', '0.9998'),
    @('E6', '  -0.03%  '),
    @('D7', '0.4462'),
    @('E7', '  -0.06%  '),
    @('D8', '0.3756'),
    @('E8', '  +2.27%  '),
    @('D9', '0.07402'),
    @('E9', '  -0.01%  '),
    @('D10', '0.8779'),
    @('E10', '  +2.60%  '),
    @('D11', '20.82'),
    @('E11', '  +0.67%  '),
    @('D12', '1.827.48'),
    @('E12', '  +1.29%  '),
    @('D13', '6.707'),
    @('E13', '  +1.45%  '),
    @('E14', '  +2.09%  '),
    @('D15', '92.83'),
    @('E15', '  +0.31%  '),
    @('D16', '0.07076'),
    @('E16', '  +0.02%  '),
    @('D17', '1.000'),
    @('E17', '  -0.12%  '),
    @('D18', '0.000008818'),
    @('E18', '  +0.89%  '),
    @('E19', '  +0.02%  '),
    @('D20', '15.08'),
    @('E20', '  +1.34%  '),
    @('D21', '27.410.76'),
    @('E21', '  +1.92%  '),
    @('D22', '5.336'),
    @('E22', '  +3.38%  '),
    @('E23', '  +0.68%  '),
    @('D24', '1.955'),
    @('E24', '  -1.88%  '),
    @('D25', '150.99'),
    @('E25', '  -0.56%  '),
    @('D26', '2.258'),
    @('E26', '  +3.22%  '),
    @('D27', '18.56'),
    @('E27', '  +0.40%  '),
    @('D28', '5.346'),
    @('E28', '  +2.55%  '),
    @('D29', '117.04'),
    @('E29', '  +0.45%  '),
    @('D30', '0.08897'),
    @('E30', '  +0.76%  '),
    @('D31', '0.7947'),
    @('E31', '  +5.43%  '),
    @('E32', '  +1.78%  '),
    @('D33', '4.550'),
    @('E33', '  +2.00%  '),
    @('D34', '2.953'),
    @('E34', '  +0.91%  '),
    @('D35', '0.9995'),
    @('E35', '  -0.02%  '),
    @('D36', '1.103'),
    @('E36', '  +1.21%  '),
    @('D37', '0.01977'),
    @('E37', '  +0.32%  '),
    @('D38', '0.05271'),
    @('E38', '  +1.32%  '),
    @('D39', '7.365'),
    @('E39', '  +5.14%  '),
    @('D40', '0.5337'),
    @('E40', '  +0.05%  '),
    @('B41', 'MXToken'),
    @('C41', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'),
    @('D41', '2.876'),
    @('E41', '  +0.22%  '),
    @('B42', 'RenderToken'),
    @('C42', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'),
    @('D42', '2.348'),
    @('E42', '  +18.57%  '),
    @('E43', '  +0.64%  '),
    @('D44', '8.671'),
    @('E44', '  +2.63%  '),
    @('D45', '0.5089'),
    @('E45', '  -1.72%  '),
    @('D46', '10.59'),
    @('E46', '  +0.40%  '),
    @('B47', 'PaxosStandard'),
    @('C47', 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'),
    @('D47', '1.001'),
    @('E47', '  -0.08%  '),
    @('B48', 'Quant'),
    @('C48', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'),
    @('D48', '105.24'),
    @('E48', '  -0.21%  '),
    @('B49', 'NEARProtocol'),
    @('C49', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'),
    @('D49', '1.684'),
    @('E49', '  +0.57%  '),
    @('B50', 'PaxDollar'),
    @('C50', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'),
    @('D50', '0.9992'),
    @('E50', '  +0.00%  '),
    @('B51', 'Cronos'),
    @('C51', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'),
    @('D51', '0.06384'),
    @('E51', '  +0.75%  ')
)

foreach ($pair in $updates) {
    $addr = $pair[0]
    $val = $pair[1]
    # All of these sheet cells are stored as literal text (inline strings) in
    # the workbook -- including ones that look numeric, e.g. "0.9999" or
    # "1.000" (trailing zeros matter) -- so force text typing via NumberFormat
    # "@" before assigning, otherwise the COM value-setter auto-coerces
    # numeric-looking strings into actual numbers. ClearFormats() afterwards
    # drops that temporary text format again so the cell keeps its original
    # (default) style, matching cells that were never touched.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}
